$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Teste1")
$ws2 = $wb.Worksheets.Item("Teste2")

# --- Sheet1 (Teste1) edits ---
# Update username strings
$ws1.Range("C2").Value = "zimmer12"
$ws1.Range("C3").Value = "morais12"
$ws1.Range("C4").Value = "zimmerfer"

# Remove row 5 (A5 "Acura") entirely
$ws1.Rows.Item(5).Delete()

# Set the active selection on sheet1
$ws1.Range("F10").Select()

# --- Sheet2 (Teste2) edits ---
# A1 becomes "Acura" (previously "modelo de carro")
$ws2.Range("A1").Value = "Acura"

# Remove row 2 (A2 "Acura") entirely
$ws2.Rows.Item(2).Delete()

# Set the active selection on sheet2
$ws2.Range("C2").Select()

# Re-select sheet1 as the active sheet (tab selected) to match original tabSelected state
$ws1.Select()
$ws1.Range("F10").Select()
